$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(291).Insert()
$ws.Range("A291").Value = 3
$ws.Range("B291").Value = "Femacal de La Calera"
$ws.Range("C291").Value = "Coquimbo"
$ws.Range("D291").Value = 45275
$ws.Range("E291").Value = 5
$ws.Range("F291").Value = 100112026
$ws.Range("G291").Value = "Haba"
$ws.Range("H291").Value = "Sin especificar"
$ws.Range("I291").Value = "Primera"
$ws.Range("J291").Value = 80
$ws.Range("K291").Value = 9000
$ws.Range("L291").Value = 10000
$ws.Range("M291").Value = 9438
$ws.Range("N291").Value = "$/saco 25 kilos"
$ws.Range("O291").Value = "Provincia de Petorca"
$ws.Range("P291").Value = 378
$ws.Range("Q291").Value = 25
$ws.Range("R291").Value = "Hortaliza"
